$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 70-73 have their data (columns A,B,E,F,G,H,Q,R,Y,AA) cyclically rotated
# downward by one row: new row70 = old row73, new row71 = old row70,
# new row72 = old row71, new row73 = old row72.
#
# Using Range.Copy (rather than re-assigning .Value/.Value2) preserves the
# exact underlying cell representation (numbers stay numbers, inline/shared
# strings stay strings, date-looking text like "2023-08-26" is not
# reinterpreted as a date serial number).

$cols = @("A","B","E","F","G","H","Q","R","Y","AA")
$bufferRow = 9999

foreach ($c in $cols) {
    # stash row 73's current value in a scratch cell far outside the used range
    $ws.Range("$c" + "73").Copy($ws.Range("$c" + $bufferRow))
    # shift 72 -> 73, 71 -> 72, 70 -> 71 (process high-to-low so sources
    # aren't clobbered before they're read)
    $ws.Range("$c" + "72").Copy($ws.Range("$c" + "73"))
    $ws.Range("$c" + "71").Copy($ws.Range("$c" + "72"))
    $ws.Range("$c" + "70").Copy($ws.Range("$c" + "71"))
    # move the stashed old row73 value into row 70
    $ws.Range("$c" + $bufferRow).Copy($ws.Range("$c" + "70"))
    # clean up the scratch cell
    $ws.Range("$c" + $bufferRow).ClearContents()
}
